$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so numeric-looking values
# (e.g. "504.50", "0.802") are not coerced into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "57.020.41"
$ws.Range("E2").Value = "  +0.89%  "

$ws.Range("D3").Value = "2.397.57"
$ws.Range("E3").Value = "  +1.78%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "504.50"
$ws.Range("E5").Value = "  -1.68%  "

$ws.Range("D6").Value = "131.85"
$ws.Range("E6").Value = "  +3.28%  "

$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "0.553"
$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("D9").Value = "2.411.07"
$ws.Range("E9").Value = "  +1.67%  "

$ws.Range("D10").Value = "0.0964"
$ws.Range("E10").Value = "  +0.69%  "

$ws.Range("E11").Value = "  -1.65%  "

$ws.Range("D12").Value = "0.320"
$ws.Range("E12").Value = "  +1.42%  "

$ws.Range("D13").Value = "4.56"
$ws.Range("E13").Value = "  -5.13%  "

$ws.Range("D14").Value = "2.823.74"
$ws.Range("E14").Value = "  +1.71%  "

$ws.Range("D15").Value = "56.901.29"
$ws.Range("E15").Value = "  +0.84%  "

$ws.Range("D16").Value = "21.73"
$ws.Range("E16").Value = "  +1.23%  "

$ws.Range("D17").Value = "0.0000134"
$ws.Range("E17").Value = "  +2.07%  "

$ws.Range("D18").Value = "2.420.76"
$ws.Range("E18").Value = "  +2.59%  "

$ws.Range("D19").Value = "10.19"
$ws.Range("E19").Value = "  -0.53%  "

$ws.Range("D20").Value = "309.91"
$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("E21").Value = "  -0.51%  "

$ws.Range("D22").Value = "6.34"
$ws.Range("E22").Value = "  +4.30%  "

$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("E24").Value = "  +0.33%  "

$ws.Range("D25").Value = "65.16"
$ws.Range("E25").Value = "  +0.70%  "

$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("B27").Value = "Polygon"
$ws.Range("C27").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D27").Value = "0.376"
$ws.Range("E27").Value = "  -3.48%  "

$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.152"
$ws.Range("E28").Value = "  -0.81%  "

$ws.Range("D29").Value = "7.46"
$ws.Range("E29").Value = "  +3.93%  "

$ws.Range("D30").Value = "172.54"
$ws.Range("E30").Value = "  -1.01%  "

$ws.Range("D31").Value = "0.0₃0723"
$ws.Range("E31").Value = "  +1.23%  "

$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("E33").Value = "  -0.19%  "

$ws.Range("D34").Value = "5.92"
$ws.Range("E34").Value = "  -3.46%  "

$ws.Range("E35").Value = "  +0.15%  "

$ws.Range("E36").Value = "  +0.02%  "

$ws.Range("D37").Value = "17.92"
$ws.Range("E37").Value = "  +1.68%  "

$ws.Range("E38").Value = "  +0.92%  "

$ws.Range("E39").Value = "  +2.80%  "

$ws.Range("D40").Value = "36.64"
$ws.Range("E40").Value = "  +3.39%  "

$ws.Range("D41").Value = "0.802"
$ws.Range("E41").Value = "  -0.12%  "

$ws.Range("D42").Value = "1.44"
$ws.Range("E42").Value = "  +0.86%  "

$ws.Range("D43").Value = "131.10"
$ws.Range("E43").Value = "  +8.08%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "4.97"
$ws.Range("E44").Value = "  +1.72%  "

$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "3.34"

$ws.Range("E46").Value = "  -0.51%  "

$ws.Range("D47").Value = "251.58"
$ws.Range("E47").Value = "  -0.72%  "

$ws.Range("E48").Value = "  +0.43%  "

$ws.Range("D49").Value = "0.0485"
$ws.Range("E49").Value = "  -0.40%  "

$ws.Range("D50").Value = "16.92"
$ws.Range("E50").Value = "  +1.93%  "

$ws.Range("E51").Value = "  +0.80%  "
